$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 39: new class-day date header (2023-01-05), same style as the A25 date header ---
$ws.Range("A25").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A39").Value = 44931

# --- Row 40: first time entry of the new day ---
$ws.Range("A26").Copy()
$ws.Range("A40").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A40").Value = 0.42708333333333331
$ws.Range("B40").Value = "講解 Plan //購買百雞計畫"

# --- Row 41: second time entry of the new day ---
$ws.Range("A26").Copy()
$ws.Range("A41").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A41").Value = 0.6875
$ws.Range("B41").Value = "扭蛋機"

$excel.CutCopyMode = 0

# --- Update sheet view: scrolled down, new active selection ---
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("B42").Select()
